$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.759.55'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '1.725.46'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '''240.77'
$ws.Range('E5').Value = '  -1.23%  '
$ws.Range('D6').Value = '''0.9976'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').Value = '''0.4854'
$ws.Range('E7').Value = '  -1.19%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = '1.728.18'
$ws.Range('E10').Value = '  +0.28%  '
$ws.Range('D11').Value = '''15.95'
$ws.Range('E11').Value = '  +3.22%  '
$ws.Range('D12').Value = '''0.06879'
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').Value = '''0.6080'
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').Value = '''4.475'
$ws.Range('E14').Value = '  -1.47%  '
$ws.Range('D15').Value = '''76.93'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').Value = '''0.9979'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '26.561.51'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D19').Value = '''0.000007149'
$ws.Range('E19').Value = '  -0.91%  '
$ws.Range('D20').Value = '''11.43'
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('D21').Value = '1.951.08'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('E22').Value = '  -1.03%  '
$ws.Range('D23').Value = '''8.582'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '''5.086'
$ws.Range('E24').Value = '  -1.42%  '
$ws.Range('D25').Value = '''137.67'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = '''15.24'
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('D27').Value = '''1.775'
$ws.Range('E27').Value = '  +3.16%  '
$ws.Range('D28').Value = '''105.96'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').Value = '''1.374'
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('D30').Value = '''3.993'
$ws.Range('E30').Value = '  +1.16%  '
$ws.Range('D31').Value = '''0.07939'
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('D32').Value = '''3.690'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').Value = '''0.04472'
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''2.596'
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '''1.007'
$ws.Range('E35').Value = '  +0.72%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''0.6208'
$ws.Range('E36').Value = '  -0.84%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '''0.9266'
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '''2.027'
$ws.Range('E38').Value = '  +3.78%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '''2.443'
$ws.Range('E39').Value = '  +2.18%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').Value = '''0.9970'
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.01494'
$ws.Range('E41').Value = '  +0.86%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''5.652'
$ws.Range('E42').Value = '  +6.88%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '''99.75'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '''0.3843'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '''6.854'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '''0.1157'
$ws.Range('E46').Value = '  -1.12%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '''0.05383'
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''7.902'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '''30.09'
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '''51.55'
$ws.Range('E50').Value = '  +1.35%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '''1.234'
$ws.Range('E51').Value = '  -0.01%  '
